$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B4 ("ratio_threshold_range" / Min) changes from 1 to 0.9
$ws.Range("B4").Value = 0.9

# Active selection moves from C3 to B4
$ws.Range("B4").Select()
